$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 29.01.2022 09:30"

# Row 8 (Benzina Albert Modrice): update Delta Cena (D8) to numeric value
$ws.Range("D8").Value = 0.2

# Row 8: update Old Datum (E8) to a numeric Excel date serial, matching the
# formatting used by the other rows in that column.
$ws.Range("E8").Value = 44590.38565972223
$ws.Range("E8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
